# Add a new "SnapGenePlasmidSource" worksheet right after "BenchlingUrlSource",
# matching the structure/columns used by the other repository-id based
# source sheets (e.g. BenchlingUrlSource / RepositoryIdSource):
#   repository_id | repository_name | input | output | type | output_name | id
# with a dropdown (list) data validation on the repository_name column.

$wb = $excel.ActiveWorkbook

$afterSheet = $wb.Worksheets.Item("BenchlingUrlSource")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$newSheet.Name = "SnapGenePlasmidSource"

$newSheet.Range("A1").Value = "repository_id"
$newSheet.Range("B1").Value = "repository_name"
$newSheet.Range("C1").Value = "input"
$newSheet.Range("D1").Value = "output"
$newSheet.Range("E1").Value = "type"
$newSheet.Range("F1").Value = "output_name"
$newSheet.Range("G1").Value = "id"

$validationRange = $newSheet.Range("B2:B1048576")
$validationRange.Validation.Add(3, 1, 1, '"addgene,genbank,benchling"')
